$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '62.730.71'
$ws.Range("E2").Value = '  +3.36%  '

# Row 3
$ws.Range("D3").Value = '2.445.99'
$ws.Range("E3").Value = '  +2.09%  '

# Row 5
$ws.Range("D5").Value = '''576.16'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.35%  '

# Row 6
$ws.Range("D6").Value = '''145.70'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.00%  '

# Row 7
$ws.Range("E7").Value = '  +0.04%  '

# Row 8
$ws.Range("E8").Value = '  +0.06%  '

# Row 9
$ws.Range("D9").Value = '2.444.61'
$ws.Range("E9").Value = '  +1.79%  '

# Row 10
$ws.Range("E10").Value = '  +3.03%  '

# Row 11
$ws.Range("E11").Value = '  +2.56%  '

# Row 12
$ws.Range("E12").Value = '  +1.20%  '

# Row 13
$ws.Range("D13").Value = '''0.353'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.23%  '

# Row 14
$ws.Range("D14").Value = '''28.14'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.94%  '

# Row 15
$ws.Range("E15").Value = '  +6.08%  '

# Row 16
$ws.Range("D16").Value = '2.888.82'
$ws.Range("E16").Value = '  +2.05%  '

# Row 17
$ws.Range("D17").Value = '62.701.13'
$ws.Range("E17").Value = '  +3.41%  '

# Row 18
$ws.Range("D18").Value = '2.446.34'
$ws.Range("E18").Value = '  +1.72%  '

# Row 19
$ws.Range("D19").Value = '''7.90'
$ws.Range("D19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = '''10.98'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.84%  '

# Row 21
$ws.Range("D21").Value = '''329.35'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.64%  '

# Row 22
$ws.Range("D22").Value = '''4.14'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.15%  '

# Row 23
$ws.Range("E23").Value = '  +9.00%  '

# Row 24
$ws.Range("E24").Value = '  +0.02%  '

# Row 25
$ws.Range("D25").Value = '''66.12'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.84%  '

# Row 26
$ws.Range("D26").Value = '''648.97'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +11.48%  '

# Row 27
$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").Value = '''1.18'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +17.71%  '

# Row 28
$ws.Range("B28").Value = 'BabyDogeCoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D28").Value = '0.0₆0580'
$ws.Range("E28").Value = '  +108.99%  '

# Row 29
$ws.Range("D29").Value = '''8.54'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.86%  '

# Row 30
$ws.Range("D30").Value = '0.0₃0991'
$ws.Range("E30").Value = '  +5.89%  '

# Row 31
$ws.Range("D31").Value = '2.571.59'

# Row 32
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").Value = '''8.20'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.13%  '

# Row 33
$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").Value = '''1.45'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +7.86%  '

# Row 34
$ws.Range("E34").Value = '  +3.51%  '

# Row 35
$ws.Range("D35").Value = '''0.138'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.84%  '

# Row 36
$ws.Range("E36").Value = '  +2.64%  '

# Row 37
$ws.Range("E37").Value = '  +0.11%  '

# Row 38
$ws.Range("E38").Value = '  +3.57%  '

# Row 39
$ws.Range("E39").Value = '  +6.97%  '

# Row 40
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").Value = '''153.58'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.70%  '

# Row 41
$ws.Range("B41").Value = 'PolygonEcosystemToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D41").Value = '''0.374'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.86%  '

# Row 42
$ws.Range("D42").Value = '''18.75'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.54%  '

# Row 43
$ws.Range("E43").Value = '  +8.31%  '

# Row 44
$ws.Range("E44").Value = '  +4.91%  '

# Row 45
$ws.Range("D45").Value = '''42.39'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.95%  '

# Row 46
$ws.Range("E46").Value = '  +0.02%  '

# Row 47
$ws.Range("E47").Value = '  +27.41%  '

# Row 48
$ws.Range("D48").Value = '''145.47'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.05%  '

# Row 49
$ws.Range("D49").Value = '''3.63'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.38%  '

# Row 50
$ws.Range("D50").Value = '''20.72'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.93%  '

# Row 51
$ws.Range("E51").Value = '  +2.73%  '

